$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions
$ws.Cells.Item(1, 12).Value = "DNA_peptidase"
$ws.Cells.Item(1, 13).Value = "DNA_sect_pep"

# Data rows: column L = DNA_peptidase value, column M = DNA_sect_pep value
$values = @{
    2 = @("3.8437106589524601E-2", "8.3316657947394099E-3")
    3 = @("4.0307176274845799E-2", "8.9464381510720205E-3")
    4 = @("3.86280959751455E-2", "1.0161982549717801E-2")
    5 = @("2.6856325675505499E-2", "4.07066236722416E-3")
    6 = @("3.5460707952912701E-2", "7.6395393154068403E-3")
    7 = @("3.5669821242514499E-2", "7.0503647795832399E-3")
    8 = @("3.0258444295056801E-2", "4.71794941396971E-3")
    9 = @("4.6222620717537799E-2", "1.1863482705294401E-2")
    10 = @("3.6934049757349499E-2", "8.9132156380819592E-3")
    11 = @("3.4733361694088798E-2", "7.7387081043315499E-3")
    12 = @("3.3956951796825799E-2", "7.8670666006048495E-3")
    13 = @("2.7726508246430401E-2", "6.5399039937055698E-3")
    14 = @("3.6335077104417898E-2", "8.6882387340151603E-3")
    15 = @("3.6053147546787498E-2", "8.7915765806195392E-3")
    16 = @("3.0597610695263198E-2", "5.62977602469035E-3")
    17 = @("2.9842164582482001E-2", "5.1658595473017604E-3")
    18 = @("3.1794812075955498E-2", "6.1266955701898898E-3")
    19 = @("2.9014651171705E-2", "4.9387261221818099E-3")
    20 = @("3.3932220939926198E-2", "6.3845509539485297E-3")
    21 = @("2.8095933368715201E-2", "6.2543950155829397E-3")
    22 = @("3.02683440965681E-2", "6.1384882361630397E-3")
    23 = @("2.9744341077088201E-2", "5.9243859739518999E-3")
    24 = @("2.5347500657914301E-2", "4.3165348522356001E-3")
    25 = @("3.0159588081919701E-2", "4.7917131072330102E-3")
    26 = @("3.0786830227507798E-2", "6.0511584784876697E-3")
    27 = @("3.1847934926133303E-2", "4.9451015046818002E-3")
    28 = @("3.3978851293892699E-2", "6.8458598250570497E-3")
    29 = @("3.0853474418609601E-2", "6.2598929847232504E-3")
    30 = @("3.39216073244097E-2", "5.4750877401099299E-3")
    31 = @("2.8496185124864501E-2", "4.3465749319269702E-3")
    32 = @("2.95478537427563E-2", "6.3119318926649296E-3")
    33 = @("4.1764251299440502E-2", "1.5555316589937701E-2")
    34 = @("3.2158222500444797E-2", "1.04419415935627E-2")
    35 = @("3.3908428310322203E-2", "6.6774787097068104E-3")
    36 = @("2.9469783222085801E-2", "6.3047809094146099E-3")
    37 = @("2.84160746165028E-2", "6.3141709139210699E-3")
    38 = @("1.96231794457381E-2", "3.5195232485872099E-3")
    39 = @("2.9540786334944599E-2", "7.03479020745576E-3")
    40 = @("3.7065001413193098E-2", "1.2473746440413901E-2")
    41 = @("1.83054814645702E-2", "4.37609649340786E-3")
    42 = @("2.98627774941643E-2", "6.5272575138982797E-3")
    43 = @("2.8685633982918599E-2", "6.1776161216754999E-3")
    44 = @("2.89751438505873E-2", "6.2910285928486197E-3")
    45 = @("2.83992838671039E-2", "6.8383826447054501E-3")
    46 = @("2.6024908790981301E-2", "6.38954000284697E-3")
    47 = @("1.6144798350799301E-2", "3.0665442025675999E-3")
    48 = @("1.78521622338264E-2", "4.6163597861979301E-3")
    49 = @("2.4629207219047499E-2", "4.1314424134086397E-3")
    50 = @("2.7933876920376199E-2", "5.8690046088717297E-3")
    51 = @("3.3486785841008997E-2", "8.1393425725092595E-3")
    52 = @("3.3537190199578903E-2", "7.7256792874708302E-3")
    53 = @("3.2066030630954402E-2", "7.5502597619708E-3")
    54 = @("3.5544252718790803E-2", "9.6844326430508398E-3")
    55 = @("3.4342369352365398E-2", "9.4744513388875899E-3")
    56 = @("3.6722480181342998E-2", "9.1698474937302096E-3")
    57 = @("4.0600503761446798E-2", "9.9481730043812708E-3")
    58 = @("3.2002647212152097E-2", "8.3014119167995008E-3")
    59 = @("2.9240910768473599E-2", "7.9534159555262902E-3")
    60 = @("3.3728807972256102E-2", "7.7262607697775997E-3")
    61 = @("3.4407030557890399E-2", "9.2137205279519806E-3")
    62 = @("2.5728734074101301E-2", "4.5717760884355198E-3")
    63 = @("2.63079551347059E-2", "5.6204617452887996E-3")
    64 = @("3.7673123570132298E-2", "1.15313739044997E-2")
    65 = @("3.3698567403086799E-2", "1.06631082392396E-2")
    66 = @("3.0297047288060701E-2", "9.4932282950734095E-3")
    67 = @("3.51094099195218E-2", "1.15643587700973E-2")
    68 = @("2.9397914443190801E-2", "8.7580355274833894E-3")
    69 = @("2.8063991257634899E-2", "7.1377456463740802E-3")
    70 = @("2.99660518447958E-2", "5.690381320548E-3")
    71 = @("3.2654823496197798E-2", "7.7993521248993902E-3")
    72 = @("2.9021987484911499E-2", "6.2870650660376099E-3")
    73 = @("3.6128081521538999E-2", "8.4433694045829392E-3")
    74 = @("3.2652523146291301E-2", "7.72936931676051E-3")
    75 = @("3.2558768474367901E-2", "8.8810103066433391E-3")
    76 = @("2.62980716899151E-2", "5.0416144563517999E-3")
    77 = @("2.5900104444965499E-2", "6.0274657987651201E-3")
    78 = @("3.2985515236838403E-2", "9.8310990853929595E-3")
    79 = @("3.3322855876399902E-2", "8.0297448847907894E-3")
    80 = @("3.9632229881838801E-2", "1.0543374555236599E-2")
    81 = @("3.9731778149778102E-2", "1.0928790470282701E-2")
    82 = @("3.4948672703872301E-2", "9.2059805304147008E-3")
    83 = @("3.2443248233743803E-2", "6.2231129180221499E-3")
    84 = @("2.7050986337532501E-2", "5.4385227612507303E-3")
    85 = @("2.75917966340926E-2", "5.04918426137163E-3")
    86 = @("2.7862350194692E-2", "5.7148754998064301E-3")
    87 = @("3.9515202561713203E-2", "8.7894288386736602E-3")
    88 = @("4.4042614066619898E-2", "1.47702920729193E-2")
    89 = @("3.7916221935526899E-2", "1.10533080385079E-2")
    90 = @("3.2391900907091403E-2", "7.93242884639077E-3")
    91 = @("3.3516245707581403E-2", "7.3728334303028199E-3")
    92 = @("2.8555125635174899E-2", "6.2935431595858302E-3")
    93 = @("3.0157076081159299E-2", "6.6194841726207702E-3")
    94 = @("3.5533759628648201E-2", "8.7542524887552504E-3")
    95 = @("3.4987440477042199E-2", "9.2661777797429599E-3")
    96 = @("3.6594592755358202E-2", "9.4499788481472499E-3")
    97 = @("3.8836136600336302E-2", "9.6080612304736603E-3")
    98 = @("3.5057462129368701E-2", "8.7376318644580894E-3")
    99 = @("3.2571601561914403E-2", "7.8134526607806093E-3")
    100 = @("3.1280241100798703E-2", "7.6834950778595099E-3")
    101 = @("2.6372683872885199E-2", "5.7939299953876297E-3")
    102 = @("3.54828569946984E-2", "8.9240155548230907E-3")
    103 = @("3.3339381210791497E-2", "8.49814798105877E-3")
    104 = @("2.8661045547340801E-2", "6.4776982983443601E-3")
    105 = @("2.8743885133800402E-2", "6.7519622006340903E-3")
    106 = @("2.99012106987964E-2", "5.4266757521490199E-3")
    107 = @("2.60428159030466E-2", "6.2490201153130899E-3")
    108 = @("3.5130273590420799E-2", "7.72243343557157E-3")
    109 = @("3.5769852748161197E-2", "8.3937893912190698E-3")
    110 = @("3.3782610376235103E-2", "7.9944875906187395E-3")
    111 = @("3.4982283913107003E-2", "8.1164814078730305E-3")
    112 = @("3.4845971724736101E-2", "8.1101962815039296E-3")
    113 = @("2.8774234398431402E-2", "6.9821857177816398E-3")
    114 = @("3.3824695972597803E-2", "8.2670858997872295E-3")
    115 = @("2.7262327495110199E-2", "6.0285061309382501E-3")
    116 = @("3.3185983137452801E-2", "7.3044518391520499E-3")
    117 = @("3.3640364139289398E-2", "8.2221451451985495E-3")
    118 = @("3.2507261823641498E-2", "6.9212937627117204E-3")
    119 = @("3.8478598242308398E-2", "1.05709167174971E-2")
    120 = @("3.6915111617383603E-2", "1.12613861580386E-2")
    121 = @("2.8588810338575299E-2", "5.2667692066990001E-3")
    122 = @("2.5553936878511999E-2", "5.2280091339347501E-3")
    123 = @("3.5437162877555897E-2", "1.09402223908653E-2")
    124 = @("3.1564738059720099E-2", "6.2384549139398501E-3")
    125 = @("3.3153262394410903E-2", "9.8227395934595493E-3")
    126 = @("3.7145777713616997E-2", "8.8720614402111805E-3")
    127 = @("3.7665839747317E-2", "9.5260431809676901E-3")
    128 = @("4.2586291150179402E-2", "1.19353957522672E-2")
    129 = @("3.9695847709189797E-2", "1.02000156358112E-2")
    130 = @("3.0341447375166E-2", "6.8603038803077996E-3")
    131 = @("4.0253422774438399E-2", "1.4598654292125699E-2")
    132 = @("4.0943635040340397E-2", "1.39807558038716E-2")
    133 = @("3.7931089665537703E-2", "1.39885037765143E-2")
    134 = @("4.9260780154097399E-2", "1.82387020375008E-2")
    135 = @("4.8353756337092997E-2", "1.7327482697513601E-2")
    136 = @("4.5575213742406202E-2", "1.65237161682181E-2")
    137 = @("5.2466800898783701E-2", "2.0650132599920198E-2")
    138 = @("3.6102929913419603E-2", "9.6970194850892998E-3")
    139 = @("2.9547456709924001E-2", "7.4568448404734302E-3")
    140 = @("4.0465235478028001E-2", "1.40341191162254E-2")
    141 = @("4.1862867656919898E-2", "1.5664474110533198E-2")
    142 = @("4.4777482052587601E-2", "1.52325524023231E-2")
    143 = @("4.3148770331021803E-2", "1.4166716430653099E-2")
    144 = @("4.2807662103137502E-2", "1.47733623732949E-2")
    145 = @("3.7489324747111502E-2", "1.1739645819578299E-2")
    146 = @("3.3592774833807397E-2", "8.0662636031163405E-3")
    147 = @("2.2651684140328299E-2", "7.0679903360454203E-3")
    148 = @("3.6518393296612398E-2", "1.1520382948812E-2")
    149 = @("3.23301988759081E-2", "7.98711492557104E-3")
    150 = @("3.5973173016671599E-2", "9.9012503723037194E-3")
    151 = @("4.0793548906474797E-2", "1.3448895938685E-2")
    152 = @("3.11155907835362E-2", "7.8343306649127607E-3")
    153 = @("3.5894951019381602E-2", "1.10635147216706E-2")
    154 = @("3.0231887017841599E-2", "7.7555969253345902E-3")
    155 = @("3.7331858297558397E-2", "1.19058825252348E-2")
    156 = @("3.2369837663803497E-2", "1.00800702361761E-2")
    157 = @("5.12765418952668E-2", "1.9760212494082201E-2")
    158 = @("3.3005921012850097E-2", "9.2533348505823496E-3")
    159 = @("3.9464306774498498E-2", "1.4411892410494599E-2")
    160 = @("4.4232529161969701E-2", "1.6398058230224499E-2")
    161 = @("3.2804327992051498E-2", "9.5259455638531595E-3")
    162 = @("4.1121587248154401E-2", "1.3596225162397699E-2")
    163 = @("3.5292451126234699E-2", "1.13514623631013E-2")
    164 = @("3.5317916499603599E-2", "1.1711558875188099E-2")
    165 = @("3.8406685823549801E-2", "1.3283871923889501E-2")
    166 = @("3.2982127376997197E-2", "1.20324401670524E-2")
    167 = @("2.80258454026991E-2", "6.7335099041055798E-3")
    168 = @("4.0735096220979303E-2", "1.31857154386408E-2")
    169 = @("3.8936051032613099E-2", "8.4327036613489198E-3")
    170 = @("3.9491596532109097E-2", "8.8349023173879392E-3")
    171 = @("3.2680618392519699E-2", "6.1817277984886404E-3")
    172 = @("3.2309973472867799E-2", "6.1323703861595504E-3")
    173 = @("3.8997539913207202E-2", "8.8719468181259198E-3")
    174 = @("3.8886410481556402E-2", "9.4091990213327498E-3")
    175 = @("3.6753164063223001E-2", "9.2628101843995492E-3")
    176 = @("3.8846271252384601E-2", "9.2322860145774201E-3")
    177 = @("2.9775079237966601E-2", "5.7477622862456799E-3")
    178 = @("4.2119029169070801E-2", "1.39919263526817E-2")
    179 = @("2.8806487064957701E-2", "5.4215437067842102E-3")
    180 = @("4.0810800594331502E-2", "1.2507861200772101E-2")
    181 = @("3.6420763067647403E-2", "7.4310540377019096E-3")
}

foreach ($r in $values.Keys) {
    $pair = $values[$r]
    $ws.Cells.Item([int]$r, 12).Value = [double]$pair[0]
    $ws.Cells.Item([int]$r, 13).Value = [double]$pair[1]
}

# Update selection to K16 as in the diff
$null = $ws.Range("K16").Select()